$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.048.89'
$ws.Range("E2").Value = '  -3.12%  '
$ws.Range("D3").Value = '2.536.94'
$ws.Range("E3").Value = '  -4.38%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '512.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.555'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.98%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.46'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.07%  '
$ws.Range("E10").Value = '  -3.25%  '
$ws.Range("E11").Value = '  -3.11%  '
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Value = '2.983.40'
$ws.Range("E13").Value = '  -4.33%  '
$ws.Range("D14").Value = '57.019.54'
$ws.Range("E14").Value = '  -3.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.97'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.11%  '
$ws.Range("E16").Value = '  -2.95%  '
$ws.Range("D17").Value = '2.530.12'
$ws.Range("E17").Value = '  -4.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '331.75'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.27'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.08'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.73%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.15'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.166'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.400'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.10%  '
$ws.Range("D27").Value = '2.659.41'
$ws.Range("E27").Value = '  -4.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.83%  '
$ws.Range("E29").Value = '  -6.67%  '
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.27'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.14%  '
$ws.Range("E32").Value = '  -2.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '148.76'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.39'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.96'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.34%  '
$ws.Range("E36").Value = '  -5.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.842'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '35.68'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.95%  '
$ws.Range("E39").Value = '  -5.29%  '
$ws.Range("E40").Value = '  -4.49%  '
$ws.Range("E41").Value = '  +0.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.46'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.35%  '
$ws.Range("E43").Value = '  -1.94%  '
$ws.Range("E44").Value = '  -0.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.574'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.90%  '
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0520'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.22%  '
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '257.54'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.60%  '
$ws.Range("D49").Value = '1.966.76'
$ws.Range("E49").Value = '  -3.91%  '
$ws.Range("E50").Value = '  -2.75%  '
$ws.Range("E51").Value = '  -5.06%  '